$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the old column T ("ID GA"), pushing it out to column W.
# This naturally carries over formatting/width and shifts everything correctly.
$ws.Range("T1:V2").EntireColumn.Insert()

# Re-extend the PayMe header merge from R1:S1 to the new R1:V1.
$ws.Range("R1:V1").Merge()

# New header labels for the 3 inserted columns (row 2).
$ws.Range("T2").Value = "Adquirier ID"
$ws.Range("U2").Value = "Wallet Password"
$ws.Range("V2").Value = "Gateway Password"

# Highlight the "Legal" block data row (A2:L2) with a yellow fill.
$ws.Range("A2:L2").Interior.Color = 65535
